# Letter-Frequency.xlsx commit: "Update letter frequency values * 100 in source document"
#
# The underlying edit re-sorts the "Data" table from a descending sort on the
# "Median" column (F) to an ascending sort on the "Letter" column (A), and
# switches the number format used by the calculated columns (Text Value /
# Dictionary Value / Median) from "0.00" to "0.0000". The active selection
# also ends up on E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-sort the table A2:F27 by column A (Letter), ascending, header row
#        already present (rows 1 stays as header). This mirrors the table's
#        sortState moving from `descending on F1:F27` to `ascending on A1:A27`.
$sortRange = $ws.Range("A1:F27")
$keyRange  = $ws.Range("A1")
$null = $sortRange.Sort($keyRange, 1, $null, $null, 2, $null, 2, 1)

# --- 2. Update the number format of the calculated columns (D:F) from
#        "0.00" (numFmtId 2) to "0.0000" (custom numFmtId 164), matching the
#        header cells too so the whole D1:F27 block is consistent.
$ws.Range("D1:F27").NumberFormat = "0.0000"

# --- 3. Leave the selection on E3, matching the saved sheet view.
$null = $ws.Range("E3").Select()
